$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B192 and B193 in place (C192/C193 stay the same).
$ws.Range("B192").Value = "5582/2"
$ws.Range("B193").Value = "9213/12"

# Rows 194-283 held the remaining "particelle" records. The last 26 of those
# records (formerly rows 220-283) move up to occupy rows 194-257, while the
# row index column (A) keeps its original sequential values. Copy just the
# B:C values up by 26 rows, then remove the now-empty trailing rows.
# Column B holds particelle codes that must stay text (some look like plain
# numbers, e.g. "45"), so force text formatting before the paste and restore
# the original (unstyled) look afterwards.
$ws.Range("B194:B257").NumberFormat = "@"
$srcValues = $ws.Range("B220:C283").Value()
$ws.Range("B194:C257").Value = $srcValues
$ws.Range("B194:B257").Style = "Normal"

# Remove the now-obsolete trailing rows 258:283 so the used range / dimension
# shrinks back down to row 257. (-4162 = xlShiftUp)
$ws.Range("A258:C283").Delete(-4162)
